# Add 2022-Q4 data
#
# The workbook has a "总计" (summary) sheet followed by one sheet per
# quarter (newest first). This change inserts a brand-new "2022-Q4"
# quarter sheet (cloned from the "2022-Q3" sheet so it keeps identical
# formatting/styles) right after "总计", fills it with the new quarter's
# two fund rows, and updates the summary sheet so its per-quarter table
# now also includes the new quarter (shifting the previously-listed
# quarters down by one row and appending the oldest quarter, 2020-Q4,
# as a new last row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by cloning "2022-Q3" (sheet index 2)
#    so that all styles/column widths/number formats carry over exactly.
#    The clone is placed immediately before "2022-Q3", i.e. right after
#    "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Make sure the fund-code / numeric-looking text columns stay text
# (they must keep any leading zeros and not turn into real numbers).
$q4.Range("B2:G3").NumberFormat = "@"

# Drop the rows that don't apply to 2022-Q4 (originally rows 4-6).
$q4.Range("A4:H6").EntireRow.Delete()

# Row 2: 004497 / 前海开源多元策略灵活配置混合C keeps its code & name,
# only the figures change.
$q4.Range("D2").Value = "1.79"
$q4.Range("E2").Value = "79.66"
$q4.Range("F2").Value = "3.45"
$q4.Range("G2").Value = "0.0618"
$q4.Range("H2").Value = 10

# Row 3: now 004496 / 前海开源多元策略灵活配置混合A (was 160135 / 南方...).
$q4.Range("B3").Value = "004496"
$q4.Range("C3").Value = "前海开源多元策略灵活配置混合A"
$q4.Range("D3").Value = "1.30"
$q4.Range("E3").Value = "79.66"
$q4.Range("F3").Value = "3.45"
$q4.Range("G3").Value = "0.0448"
$q4.Range("H3").Value = 10

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: shift the existing quarter rows
#    down by one (the newest data, 2022-Q4, becomes the new row 2) and
#    append the now-overflowed last quarter (2020-Q4) as a new row 10.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.11

$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 5
$summary.Range("D3").Value = 0.17

$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 5
$summary.Range("D4").Value = 0.29

$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 0.23

$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 6
$summary.Range("D6").Value = 0.31

$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 5
$summary.Range("D7").Value = 0.11

$summary.Range("B8").Value = "2021-Q2"
$summary.Range("C8").Value = 4
$summary.Range("D8").Value = 0.26

$summary.Range("B9").Value = "2021-Q1"
$summary.Range("C9").Value = 3
$summary.Range("D9").Value = 0.2

# New row 10 - carry the formatting from row 9 before filling it in.
$summary.Range("A9:D9").Copy()
$summary.Range("A10:D10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$summary.Range("A10").Value = 8
$summary.Range("B10").Value = "2020-Q4"
$summary.Range("C10").Value = 3
$summary.Range("D10").Value = 0.19

# ---------------------------------------------------------------------
# 3) Keep the originally-selected tab (the oldest quarter, now the last
#    sheet) as the active sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
